$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.808940836089925
$ws.Range("C2").Value = 6.336820006395645
$ws.Range("E2").Value = 13.16751500615235
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 23.46375690846661
$ws.Range("H2").Value = 13.07585118885329
$ws.Range("K2").Value = 8.15568647088423
$ws.Range("M2").Value = 12.9348255823596
$ws.Range("N2").Value = 17.79635808369551
$ws.Range("O2").Value = 19.17627809140492

$ws.Range("B3").Value = 7.517273681607985
$ws.Range("C3").Value = 6.279838658053961
$ws.Range("E3").Value = 12.9515914680131
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 23.55791194171996
$ws.Range("H3").Value = 13.12165793980789
$ws.Range("K3").Value = 7.925188406367674
$ws.Range("M3").Value = 12.75956877088353
$ws.Range("N3").Value = 17.84937946615135
$ws.Range("O3").Value = 19.25517580693503

$ws.Range("B4").Value = 7.333251365390312
$ws.Range("C4").Value = 6.244450165750997
$ws.Range("E4").Value = 12.82179424058236
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 23.62355177637745
$ws.Range("H4").Value = 13.15167176831013
$ws.Range("K4").Value = 7.778730925400929
$ws.Range("M4").Value = 12.65370863413821
$ws.Range("N4").Value = 17.88356354045589
$ws.Range("O4").Value = 19.30750560318572

$ws.Range("B5").Value = 7.257134457001943
$ws.Range("C5").Value = 6.229936237903948
$ws.Range("E5").Value = 12.76967705153642
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 23.65225902559355
$ws.Range("H5").Value = 13.1643778074053
$ws.Range("K5").Value = 7.717861435591747
$ws.Range("M5").Value = 12.61106198792493
$ws.Range("N5").Value = 17.89790447352193
$ws.Range("O5").Value = 19.32980644482137

$ws.Range("B6").Value = 7.244430879186946
$ws.Range("C6").Value = 6.227520844287315
$ws.Range("E6").Value = 12.76107213782246
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 23.65714382984622
$ws.Range("H6").Value = 13.16651634402385
$ws.Range("K6").Value = 7.70768405570055
$ws.Range("M6").Value = 12.60401179479384
$ws.Range("N6").Value = 17.90031060857819
$ws.Range("O6").Value = 19.33356839682995

$ws.Range("B7").Value = 7.332229225463102
$ws.Range("C7").Value = 6.244254790872892
$ws.Range("E7").Value = 12.82108812626282
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 23.62393101607719
$ws.Range("H7").Value = 13.15184120204416
$ws.Range("K7").Value = 7.77791475047943
$ws.Range("M7").Value = 12.65313142468293
$ws.Range("N7").Value = 17.88375528316253
$ws.Range("O7").Value = 19.30780240984821

$ws.Range("B8").Value = 7.709464243182443
$ws.Range("C8").Value = 6.317261107998384
$ws.Range("E8").Value = 13.09253400159935
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 23.49459092171411
$ws.Range("H8").Value = 13.09125363946116
$ws.Range("K8").Value = 8.077264516598566
$ws.Range("M8").Value = 12.87406596536883
$ws.Range("N8").Value = 17.81430246938883
$ws.Range("O8").Value = 19.20267469883587

$ws.Range("B9").Value = 8.405577227446461
$ws.Range("C9").Value = 6.45688341522709
$ws.Range("E9").Value = 13.64347326951698
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 23.30349696532243
$ws.Range("H9").Value = 12.98740936412596
$ws.Range("K9").Value = 8.623065051803716
$ws.Range("M9").Value = 13.31885022681607
$ws.Range("N9").Value = 17.69097884726695
$ws.Range("O9").Value = 19.02740954051824

$ws.Range("B10").Value = 8.885212895353581
$ws.Range("C10").Value = 6.556835264816828
$ws.Range("E10").Value = 14.05501425397536
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 23.20178232651834
$ws.Range("H10").Value = 12.92021884041896
$ws.Range("K10").Value = 8.996429596846903
$ws.Range("M10").Value = 13.64952229507914
$ws.Range("N10").Value = 17.60814994982367
$ws.Range("O10").Value = 18.91754741814109

$ws.Range("B11").Value = 9.095617801709366
$ws.Range("C11").Value = 6.601637674087748
$ws.Range("E11").Value = 14.24274889712503
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 23.16401717413054
$ws.Range("H11").Value = 12.89162477804802
$ws.Range("K11").Value = 9.159835165481258
$ws.Range("M11").Value = 13.80013963089687
$ws.Range("N11").Value = 17.57214262658616
$ws.Range("O11").Value = 18.87168942979028

$ws.Range("B12").Value = 9.174111713921077
$ws.Range("C12").Value = 6.618499322574528
$ws.Range("E12").Value = 14.31383726712016
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 23.15094750955027
$ws.Range("H12").Value = 12.8810801589676
$ws.Range("K12").Value = 9.220754258365025
$ws.Range("M12").Value = 13.85714854001564
$ws.Range("N12").Value = 17.55874691530119
$ws.Range("O12").Value = 18.85491789417624

$ws.Range("B13").Value = 9.157260117295996
$ws.Range("C13").Value = 6.614872622222457
$ws.Range("E13").Value = 14.29852844321136
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 23.15370741625164
$ws.Range("H13").Value = 12.88333853202536
$ws.Range("K13").Value = 9.207677383629004
$ws.Range("M13").Value = 13.844872684842
$ws.Range("N13").Value = 17.56162128534424
$ws.Range("O13").Value = 18.85850350549914

$ws.Range("B14").Value = 9.102099558211908
$ws.Range("C14").Value = 6.603027024081765
$ws.Range("E14").Value = 14.24859786958286
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 23.16291720608758
$ws.Range("H14").Value = 12.89075158819713
$ws.Range("K14").Value = 9.16486641314939
$ws.Range("M14").Value = 13.80483062763156
$ws.Range("N14").Value = 17.57103575917326
$ws.Range("O14").Value = 18.87029771383859

$ws.Range("B15").Value = 9.068156495350562
$ws.Range("C15").Value = 6.595757460500973
$ws.Range("E15").Value = 14.21801134103823
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 23.16871903411689
$ws.Range("H15").Value = 12.89532919367638
$ws.Range("K15").Value = 9.138517642129159
$ws.Range("M15").Value = 13.78029865754109
$ws.Range("N15").Value = 17.5768335558595
$ws.Range("O15").Value = 18.87759939656906

$ws.Range("B16").Value = 8.871300270159296
$ws.Range("C16").Value = 6.553893207190447
$ws.Range("E16").Value = 14.0427491576952
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 23.20442222616485
$ws.Range("H16").Value = 12.92212718460548
$ws.Range("K16").Value = 8.985618169011587
$ws.Range("M16").Value = 13.63967830276206
$ws.Range("N16").Value = 17.6105366760536
$ws.Range("O16").Value = 18.92062736431184

$ws.Range("B17").Value = 8.748494526842808
$ws.Range("C17").Value = 6.528034291996988
$ws.Range("E17").Value = 13.93530915135564
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 23.22850938135281
$ws.Range("H17").Value = 12.9390716693254
$ws.Range("K17").Value = 8.890145264685076
$ws.Range("M17").Value = 13.55342258820091
$ws.Range("N17").Value = 17.63164000044285
$ws.Range("O17").Value = 18.9480797044258

$ws.Range("B18").Value = 8.677130441362833
$ws.Range("C18").Value = 6.513098707031316
$ws.Range("E18").Value = 13.87356619486967
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 23.24316379313257
$ws.Range("H18").Value = 12.94900323845792
$ws.Range("K18").Value = 8.834627364170302
$ws.Range("M18").Value = 13.50383191827554
$ws.Range("N18").Value = 17.64393550264051
$ws.Range("O18").Value = 18.96425715313171

$ws.Range("B19").Value = 8.652844546117949
$ws.Range("C19").Value = 6.508031342295988
$ws.Range("E19").Value = 13.85267272612537
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 23.24826268819889
$ws.Range("H19").Value = 12.95239776968259
$ws.Range("K19").Value = 8.815727213606277
$ws.Range("M19").Value = 13.48704671251097
$ws.Range("N19").Value = 17.64812561588347
$ws.Range("O19").Value = 18.96980107918301

$ws.Range("B20").Value = 8.761643385055553
$ws.Range("C20").Value = 6.530793517616853
$ws.Range("E20").Value = 13.9467413384412
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 23.22586239256458
$ws.Range("H20").Value = 12.93724869729217
$ws.Range("K20").Value = 8.900371323064892
$ws.Range("M20").Value = 13.56260284625519
$ws.Range("N20").Value = 17.62937722919309
$ws.Range("O20").Value = 18.94511723034582

$ws.Range("B21").Value = 9.118334084335538
$ws.Range("C21").Value = 6.60650925226993
$ws.Range("E21").Value = 14.26326436840625
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 23.16017859221405
$ws.Range("H21").Value = 12.8885665053314
$ws.Range("K21").Value = 9.177467306721237
$ws.Range("M21").Value = 13.81659309862577
$ws.Range("N21").Value = 17.56826400802815
$ws.Range("O21").Value = 18.86681733817177

$ws.Range("B22").Value = 9.344539959792975
$ws.Range("C22").Value = 6.655383664330976
$ws.Range("E22").Value = 14.47008084136701
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 23.12442998740133
$ws.Range("H22").Value = 12.8584013026026
$ws.Range("K22").Value = 9.352962619729976
$ws.Range("M22").Value = 13.98241118208318
$ws.Range("N22").Value = 17.52971843241022
$ws.Range("O22").Value = 18.8191062188032

$ws.Range("B23").Value = 9.224460648155791
$ws.Range("C23").Value = 6.629356986199491
$ws.Range("E23").Value = 14.35972852470926
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 23.14285024704692
$ws.Range("H23").Value = 12.87434997305971
$ws.Range("K23").Value = 9.259820135305837
$ws.Range("M23").Value = 13.8939449227826
$ws.Range("N23").Value = 17.55016355668071
$ws.Range("O23").Value = 18.84425320210813

$ws.Range("B24").Value = 8.755701153094979
$ws.Range("C24").Value = 6.529546285863123
$ws.Range("E24").Value = 13.94157276054442
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 23.22705658498479
$ws.Range("H24").Value = 12.93807227122275
$ws.Range("K24").Value = 8.895750079546826
$ws.Range("M24").Value = 13.55845245087987
$ws.Range("N24").Value = 17.6303997205107
$ws.Range("O24").Value = 18.9464553354256

$ws.Range("B25").Value = 8.222510274468741
$ws.Range("C25").Value = 6.419539924368712
$ws.Range("E25").Value = 13.49291166708775
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 23.34843811001003
$ws.Range("H25").Value = 13.01390180475591
$ws.Range("K25").Value = 8.480099658334831
$ws.Range("M25").Value = 13.19761506146762
$ws.Range("N25").Value = 17.72297037026
$ws.Range("O25").Value = 19.07150883027865

